$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 396, shifting the existing data (old rows 396-410)
# down to become rows 398-412.
$ws.Rows("396:397").Insert()

# New row 396
$ws.Range("A396").Value = 4
$ws.Range("B396").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C396").Value = "Los Lagos"
$ws.Range("D396").Value = 45041
$ws.Range("E396").Value = 10
$ws.Range("F396").Value = "Fruta"
$ws.Range("G396").Value = 100108
$ws.Range("H396").Value = "Tropicales y subtropicales"
$ws.Range("I396").Value = 100108005
$ws.Range("J396").Value = "Piña"
$ws.Range("K396").Value = "Caramelo"
$ws.Range("L396").Value = "Primera"
$ws.Range("M396").Value = 300
$ws.Range("N396").Value = 21000
$ws.Range("O396").Value = 22000
$ws.Range("P396").Value = 21500
$ws.Range("Q396").Value = "$/caja 12 unidades"
$ws.Range("R396").Value = "Ecuador"
$ws.Range("S396").Value = 1792
$ws.Range("T396").Value = 12

# New row 397
$ws.Range("A397").Value = 4
$ws.Range("B397").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C397").Value = "Los Lagos"
$ws.Range("D397").Value = 45041
$ws.Range("E397").Value = 10
$ws.Range("F397").Value = "Fruta"
$ws.Range("G397").Value = 100108
$ws.Range("H397").Value = "Tropicales y subtropicales"
$ws.Range("I397").Value = 100108005
$ws.Range("J397").Value = "Piña"
$ws.Range("K397").Value = "Caramelo"
$ws.Range("L397").Value = "Segunda"
$ws.Range("M397").Value = 200
$ws.Range("N397").Value = 18000
$ws.Range("O397").Value = 19000
$ws.Range("P397").Value = 18500
$ws.Range("Q397").Value = "$/caja 14 unidades"
$ws.Range("R397").Value = "Ecuador"
$ws.Range("S397").Value = 1321
$ws.Range("T397").Value = 14
